# Update TPM-derived NATMI ligand-receptor statistics with recalculated values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 7.649422333333334
$ws.Range("H2").Value2 = 22.948267
$ws.Range("I2").Value2 = 0.004484559810904267
$ws.Range("J2").Value2 = 0.004484559810904268
$ws.Range("M2").Value2 = 1.343359
$ws.Range("N2").Value2 = 4.030077
$ws.Range("O2").Value2 = 0.736296379391111
$ws.Range("P2").Value2 = 0.7362963793911109
$ws.Range("Q2").Value2 = 10.27592033628433
$ws.Range("R2").Value2 = 92.48328302655901
$ws.Range("S2").Value2 = 0.003301965151931698
$ws.Range("T2").Value2 = 0.003301965151931698
$ws.Range("G3").Value2 = 7.649422333333334
$ws.Range("H3").Value2 = 22.948267
$ws.Range("I3").Value2 = 0.004484559810904267
$ws.Range("J3").Value2 = 0.004484559810904268
$ws.Range("O3").Value2 = 0.1764523396969075
$ws.Range("P3").Value2 = 0.1764523396969075
$ws.Range("Q3").Value2 = 2.462609129459334
$ws.Range("R3").Value2 = 22.163482165134
$ws.Range("S3").Value2 = 0.0007913110711447788
$ws.Range("T3").Value2 = 0.000791311071144779
$ws.Range("G4").Value2 = 7.649422333333334
$ws.Range("H4").Value2 = 22.948267
$ws.Range("I4").Value2 = 0.004484559810904267
$ws.Range("J4").Value2 = 0.004484559810904268
$ws.Range("O4").Value2 = 0.08725128091198156
$ws.Range("P4").Value2 = 0.08725128091198156
$ws.Range("Q4").Value2 = 1.217698792206111
$ws.Range("R4").Value2 = 10.959289129855
$ws.Range("S4").Value2 = 0.0003912835878277911
$ws.Range("T4").Value2 = 0.0003912835878277912
$ws.Range("I5").Value2 = 0.8893308176045429
$ws.Range("J5").Value2 = 0.889330817604543
$ws.Range("M5").Value2 = 1.343359
$ws.Range("N5").Value2 = 4.030077
$ws.Range("O5").Value2 = 0.736296379391111
$ws.Range("P5").Value2 = 0.7362963793911109
$ws.Range("Q5").Value2 = 2037.812632599089
$ws.Range("R5").Value2 = 18340.3136933918
$ws.Range("S5").Value2 = 0.6548110610831614
$ws.Range("T5").Value2 = 0.6548110610831614
$ws.Range("I6").Value2 = 0.8893308176045429
$ws.Range("J6").Value2 = 0.889330817604543
$ws.Range("O6").Value2 = 0.1764523396969075
$ws.Range("P6").Value2 = 0.1764523396969075
$ws.Range("S6").Value2 = 0.1569245035308852
$ws.Range("T6").Value2 = 0.1569245035308853
$ws.Range("I7").Value2 = 0.8893308176045429
$ws.Range("J7").Value2 = 0.889330817604543
$ws.Range("O7").Value2 = 0.08725128091198156
$ws.Range("P7").Value2 = 0.08725128091198156
$ws.Range("S7").Value2 = 0.07759525299049622
$ws.Range("T7").Value2 = 0.07759525299049622
$ws.Range("I8").Value2 = 0.1061846225845528
$ws.Range("J8").Value2 = 0.1061846225845528
$ws.Range("M8").Value2 = 1.343359
$ws.Range("N8").Value2 = 4.030077
$ws.Range("O8").Value2 = 0.736296379391111
$ws.Range("P8").Value2 = 0.7362963793911109
$ws.Range("Q8").Value2 = 243.3114438487697
$ws.Range("R8").Value2 = 2189.802994638927
$ws.Range("S8").Value2 = 0.07818335315601782
$ws.Range("T8").Value2 = 0.07818335315601782
$ws.Range("I9").Value2 = 0.1061846225845528
$ws.Range("J9").Value2 = 0.1061846225845528
$ws.Range("O9").Value2 = 0.1764523396969075
$ws.Range("P9").Value2 = 0.1764523396969075
$ws.Range("S9").Value2 = 0.01873652509487742
$ws.Range("T9").Value2 = 0.01873652509487742
$ws.Range("I10").Value2 = 0.1061846225845528
$ws.Range("J10").Value2 = 0.1061846225845528
$ws.Range("O10").Value2 = 0.08725128091198156
$ws.Range("P10").Value2 = 0.08725128091198156
$ws.Range("S10").Value2 = 0.009264744333657557
$ws.Range("T10").Value2 = 0.009264744333657558
